$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated transition-probability matrix values (re-simulated with more games)

# Row 2
$ws.Range("B2").Value = 0.1825396825396825
$ws.Range("C2").Value = 0.5952380952380952
$ws.Range("J2").Value = 0.01984126984126984
$ws.Range("P2").Value = 0.1071428571428571
$ws.Range("S2").Value = 0.09523809523809523

# Row 3
$ws.Range("C3").Value = 0.03184713375796178
$ws.Range("J3").Value = 0.04458598726114649
$ws.Range("P3").Value = 0.7388535031847133
$ws.Range("S3").Value = 0.1847133757961783

# Row 4
$ws.Range("J4").Value = 0.05128205128205128
$ws.Range("P4").Value = 0.717948717948718
$ws.Range("S4").Value = 0.2307692307692308

# Row 5
$ws.Range("P5").Value = 0.8
$ws.Range("S5").Value = 0.2

# Row 6
$ws.Range("B6").Value = 0.0410958904109589
$ws.Range("D6").Value = 0.0182648401826484
$ws.Range("F6").Value = 0.0776255707762557
$ws.Range("J6").Value = 0.2557077625570776
$ws.Range("O6").Value = 0.0136986301369863
$ws.Range("Q6").Value = 0.1735159817351598
$ws.Range("R6").Value = 0.0593607305936073
$ws.Range("S6").Value = 0.3607305936073059

# Row 7
$ws.Range("B7").Value = 0.09405940594059406
$ws.Range("D7").Value = 0.0198019801980198
$ws.Range("F7").Value = 0.07920792079207921
$ws.Range("J7").Value = 0.09900990099009901
$ws.Range("O7").Value = 0.009900990099009901
$ws.Range("Q7").Value = 0.1633663366336634
$ws.Range("R7").Value = 0.09405940594059406
$ws.Range("S7").Value = 0.4405940594059406

# Row 8
$ws.Range("B8").Value = 0.08649789029535865
$ws.Range("D8").Value = 0.008438818565400843
$ws.Range("E8").Value = 0.002109704641350211
$ws.Range("F8").Value = 0.05907172995780591
$ws.Range("J8").Value = 0.1118143459915612
$ws.Range("O8").Value = 0.01687763713080169
$ws.Range("Q8").Value = 0.2172995780590717
$ws.Range("R8").Value = 0.1181434599156118
$ws.Range("S8").Value = 0.379746835443038

# Row 9
$ws.Range("B9").Value = 0.07725321888412018
$ws.Range("D9").Value = 0.01716738197424893
$ws.Range("F9").Value = 0.0815450643776824
$ws.Range("J9").Value = 0.09442060085836911
$ws.Range("O9").Value = 0.0128755364806867
$ws.Range("Q9").Value = 0.1759656652360515
$ws.Range("R9").Value = 0.1244635193133047
$ws.Range("S9").Value = 0.4163090128755365

# Row 10
$ws.Range("B10").Value = 0.09111791730474732
$ws.Range("D10").Value = 0.01837672281776417
$ws.Range("E10").Value = 0.003062787136294028
$ws.Range("F10").Value = 0.05742725880551302
$ws.Range("J10").Value = 0.1049004594180704
$ws.Range("O10").Value = 0.008422664624808576
$ws.Range("Q10").Value = 0.2243491577335375
$ws.Range("R10").Value = 0.08728943338437979
$ws.Range("S10").Value = 0.4050535987748852

# Row 11
$ws.Range("G11").Value = 0.1644295302013423
$ws.Range("J11").Value = 0.1006711409395973
$ws.Range("K11").Value = 0.2181208053691275
$ws.Range("L11").Value = 0.5033557046979866
$ws.Range("S11").Value = 0.01342281879194631

# Row 12
$ws.Range("G12").Value = 0.7712418300653595
$ws.Range("J12").Value = 0.1699346405228758
$ws.Range("K12").Value = 0.0130718954248366
$ws.Range("L12").Value = 0.0196078431372549
$ws.Range("S12").Value = 0.0261437908496732

# Row 13
$ws.Range("G13").Value = 0.8723404255319149
$ws.Range("J13").Value = 0.1276595744680851

# Row 14
$ws.Range("G14").Value = 0.3333333333333333
$ws.Range("S14").Value = 0.6666666666666666

# Row 15
$ws.Range("F15").Value = 0.009478672985781991
$ws.Range("H15").Value = 0.1658767772511848
$ws.Range("I15").Value = 0.0947867298578199
$ws.Range("J15").Value = 0.4075829383886256
$ws.Range("K15").Value = 0.03791469194312796
$ws.Range("M15").Value = 0.04265402843601896
$ws.Range("O15").Value = 0.04739336492890995
$ws.Range("S15").Value = 0.1943127962085308

# Row 16
$ws.Range("F16").Value = 0.03571428571428571
$ws.Range("H16").Value = 0.1904761904761905
$ws.Range("I16").Value = 0.1130952380952381
$ws.Range("J16").Value = 0.4047619047619048
$ws.Range("K16").Value = 0.1011904761904762
$ws.Range("M16").Value = 0.02380952380952381
$ws.Range("N16").Value = 0.005952380952380952
$ws.Range("O16").Value = 0.04761904761904762
$ws.Range("S16").Value = 0.07738095238095238

# Row 17
$ws.Range("F17").Value = 0.01391650099403579
$ws.Range("H17").Value = 0.1848906560636183
$ws.Range("I17").Value = 0.1053677932405567
$ws.Range("J17").Value = 0.4433399602385686
$ws.Range("K17").Value = 0.07952286282306163
$ws.Range("M17").Value = 0.01391650099403579
$ws.Range("N17").Value = 0.001988071570576541
$ws.Range("O17").Value = 0.07157057654075547
$ws.Range("S17").Value = 0.08548707753479125

# Row 18
$ws.Range("F18").Value = 0.004347826086956522
$ws.Range("H18").Value = 0.1521739130434783
$ws.Range("I18").Value = 0.1043478260869565
$ws.Range("J18").Value = 0.4391304347826087
$ws.Range("K18").Value = 0.1173913043478261
$ws.Range("M18").Value = 0.02608695652173913
$ws.Range("O18").Value = 0.07391304347826087
$ws.Range("S18").Value = 0.08260869565217391

# Row 19
$ws.Range("F19").Value = 0.01671974522292994
$ws.Range("H19").Value = 0.2245222929936306
$ws.Range("I19").Value = 0.0963375796178344
$ws.Range("J19").Value = 0.3853503184713376
$ws.Range("K19").Value = 0.106687898089172
$ws.Range("M19").Value = 0.01910828025477707
$ws.Range("N19").Value = 0.0007961783439490446
$ws.Range("O19").Value = 0.06847133757961783
$ws.Range("S19").Value = 0.08200636942675159
